# Fruta / hortaliza, semanal
# A new weekly price record is inserted as row 34 (pushing the existing
# rows 34-47 down to 35-48), and the sheet's used-range grows from
# A1:R47 to A1:R48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 34..47 down to 35..48, leaving a blank row 34 to fill in.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new weekly record.
$ws.Cells.Item(34, 1).Value = 1
$ws.Cells.Item(34, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(34, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(34, 4).Value = 44900
$ws.Cells.Item(34, 5).Value = 15
$ws.Cells.Item(34, 6).Value = 100112052
$ws.Cells.Item(34, 7).Value = "Albahaca"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 350
$ws.Cells.Item(34, 11).Value = 1300
$ws.Cells.Item(34, 12).Value = 1500
$ws.Cells.Item(34, 13).Value = 1414
$ws.Cells.Item(34, 14).Value = "`$/paquete"
$ws.Cells.Item(34, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(34, 16).Value = 1414
$ws.Cells.Item(34, 17).Value = 1
$ws.Cells.Item(34, 18).Value = "Hortaliza"
